# Update crypto price/volume data per latest symbol-list refresh (Jan 15 2023 14:19 UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume 1h) are stored as text (e.g. "298.14", "-2.23%")
# rather than numbers, so force text formatting before assignment, then drop back to the
# default "Normal" style afterwards to avoid leaving a stray explicit style on the cell.
$cells = @(
    @{Addr="D2"; Value='298.14'},
    @{Addr="E2"; Value='-2.23%'},
    @{Addr="D3"; Value='31.45'},
    @{Addr="E3"; Value='-1.57%'},
    @{Addr="D4"; Value='5.147'},
    @{Addr="E4"; Value='-2.45%'},
    @{Addr="D5"; Value='0.07322'},
    @{Addr="E5"; Value='-2.34%'},
    @{Addr="D6"; Value='1.871'},
    @{Addr="E6"; Value='27.10%'},
    @{Addr="D7"; Value='7.756'},
    @{Addr="E7"; Value='-1.19%'},
    @{Addr="D8"; Value='3.742'},
    @{Addr="E8"; Value='-0.37%'},
    @{Addr="D9"; Value='0.9261'},
    @{Addr="E9"; Value='1.16%'},
    @{Addr="D10"; Value='0.1674'},
    @{Addr="E10"; Value='-1.16%'},
    @{Addr="D11"; Value='0.07151'},
    @{Addr="E11"; Value='-8.55%'},
    @{Addr="D12"; Value='0.07996'},
    @{Addr="E12"; Value='-1.19%'},
    @{Addr="D13"; Value='0.03003'},
    @{Addr="E13"; Value='-0.42%'},
    @{Addr="D14"; Value='0.09932'},
    @{Addr="E14"; Value='0.43%'},
    @{Addr="D15"; Value='0.001500'},
    @{Addr="E15"; Value='0.27%'},
    @{Addr="D16"; Value='0.006104'},
    @{Addr="E16"; Value='-1.79%'},
    @{Addr="D17"; Value='3.455'},
    @{Addr="E17"; Value='-0.71%'},
    @{Addr="D18"; Value='2.219'},
    @{Addr="E18"; Value='-0.56%'},
    @{Addr="E20"; Value='-1.89%'},
    @{Addr="E21"; Value='1.96%'},
    @{Addr="D22"; Value='0.04643'},
    @{Addr="E22"; Value='1.91%'},
    @{Addr="E23"; Value='-2.18%'},
    @{Addr="D24"; Value='0.001219'},
    @{Addr="E24"; Value='0.33%'},
    @{Addr="D25"; Value='0.004734'},
    @{Addr="E25"; Value='6.77%'},
    @{Addr="E26"; Value='-7.03%'},
    @{Addr="D27"; Value='0.0001875'},
    @{Addr="E27"; Value='7.86%'},
    @{Addr="D39"; Value='0.01717'},
    @{Addr="E39"; Value='-1.15%'},
    @{Addr="D40"; Value='0.04469'},
    @{Addr="E40"; Value='-1.13%'},
    @{Addr="D41"; Value='0.007090'},
    @{Addr="E41"; Value='-2.53%'},
    @{Addr="D42"; Value='0.1331'},
    @{Addr="E42"; Value='-1.22%'},
    @{Addr="D43"; Value='0.002179'},
    @{Addr="E43"; Value='-2.99%'},
    @{Addr="D44"; Value='0.01059'},
    @{Addr="E44"; Value='-24.26%'},
    @{Addr="D45"; Value='0.00006206'},
    @{Addr="E45"; Value='0.31%'},
    @{Addr="E46"; Value='-21.30%'},
    @{Addr="D47"; Value='0.7392'},
    @{Addr="E47"; Value='4.23%'}
)

foreach ($item in $cells) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

